$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7751099467277527
$ws.Range("B1").Value = 3.109837770462036
$ws.Range("C1").Value = 3.709166288375854
$ws.Range("D1").Value = 3.048982858657837
$ws.Range("E1").Value = 1.811389684677124
